$wb = $excel.ActiveWorkbook

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1240
$ws.Range("I98").Value = 1005
$ws.Range("J98").Value = 3433.3333
$ws.Range("K98").Value = 1005
$ws.Range("L98").Value = 3433.3333
$ws.Range("M98").Value = 493
$ws.Range("N98").Value = -6429.3333

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1240
$ws.Range("I122").Value = 1005
$ws.Range("J122").Value = 3433.3333
$ws.Range("K122").Value = 3015
$ws.Range("L122").Value = 10299.9999
$ws.Range("M122").Value = -565
$ws.Range("N122").Value = -15199.9999

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 465251.53
$ws.Range("I141").Value = 1743
$ws.Range("J141").Value = 928760.0600000001
$ws.Range("K141").Value = 5229
$ws.Range("L141").Value = 2786280.18
$ws.Range("M141").Value = -49
$ws.Range("N141").Value = -2796640.18

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2676.4
$ws.Range("I132").Value = 2547.5789
$ws.Range("J132").Value = 3084.3333
$ws.Range("K132").Value = 7642.736699999999
$ws.Range("L132").Value = 9252.999899999999
$ws.Range("M132").Value = -5112.736699999999
$ws.Range("N132").Value = -14312.9999

# Sheet CRP, row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 27417.334
$ws.Range("I4").Value = 4833.3335
$ws.Range("K4").Value = 4833.3335
$ws.Range("M4").Value = -4721.3335

# Sheet CRP, row 135
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

# Sheet CUL, row 20
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 900
$ws.Range("I20").Value = 900
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2700
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -2473
$ws.Range("N20").ClearContents()

# Sheet CUL, row 21
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 1160.7142
$ws.Range("I21").Value = 350
$ws.Range("J21").Value = 2241.6667
$ws.Range("K21").Value = 1050
$ws.Range("L21").Value = 6725.000100000001
$ws.Range("M21").Value = -877
$ws.Range("N21").Value = -7071.000100000001

# Sheet CUL, row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 29751
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 29751
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 89253
$ws.Range("N22").Value = -89591
$ws.Range("M22").ClearContents()

# Sheet CUL, row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 14755.857
$ws.Range("I26").Value = 72.75
$ws.Range("J26").Value = 34333.332
$ws.Range("K26").Value = 218.25
$ws.Range("L26").Value = 102999.996
$ws.Range("M26").Value = 69.75
$ws.Range("N26").Value = -103575.996

# Sheet CUL, row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 29751
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 29751
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 89253
$ws.Range("N27").Value = -89457
$ws.Range("M27").ClearContents()

# Sheet CUL, row 51
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1433
$ws.Range("I51").Value = 651.1667
$ws.Range("J51").Value = 2996.6667
$ws.Range("K51").Value = 1953.5001
$ws.Range("L51").Value = 8990.000100000001
$ws.Range("M51").Value = -1493.5001
$ws.Range("N51").Value = -9910.000100000001

# Sheet CUL, row 59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 2601.2
$ws.Range("J59").Value = 3001.5
$ws.Range("L59").Value = 9004.5
$ws.Range("N59").Value = -10084.5

# Sheet CUL, row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 15629.5
$ws.Range("I63").Value = 3012
$ws.Range("J63").Value = 23200
$ws.Range("K63").Value = 9036
$ws.Range("L63").Value = 69600
$ws.Range("M63").Value = -8287
$ws.Range("N63").Value = -71098

# Sheet CUL, row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 15629.5
$ws.Range("I66").Value = 3012
$ws.Range("J66").Value = 23200
$ws.Range("K66").Value = 27108
$ws.Range("L66").Value = 208800
$ws.Range("M66").Value = -23364
$ws.Range("N66").Value = -216288

# Sheet CUL, row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 2633.3333
$ws.Range("I116").Value = 900
$ws.Range("J116").Value = 3500
$ws.Range("K116").Value = 2700
$ws.Range("L116").Value = 10500
$ws.Range("M116").Value = 742
$ws.Range("N116").Value = -17384

# Sheet CUL, row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1669.875
$ws.Range("I118").Value = 344.75
$ws.Range("J118").Value = 2995
$ws.Range("K118").Value = 1034.25
$ws.Range("L118").Value = 8985
$ws.Range("M118").Value = 208.75
$ws.Range("N118").Value = -11471

# Sheet CUL, row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2414.7058
$ws.Range("I136").Value = 1781.1111
$ws.Range("J136").Value = 3127.5
$ws.Range("K136").Value = 5343.3333
$ws.Range("L136").Value = 9382.5
$ws.Range("M136").Value = -243.3333000000002
$ws.Range("N136").Value = -19582.5

# Sheet GSM, row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10517.777
$ws.Range("I5").Value = 3853.3333
$ws.Range("K5").Value = 3853.3333
$ws.Range("M5").Value = -3741.3333

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 25316.28
$ws.Range("I102").Value = 1642.7576
$ws.Range("J102").Value = 103438.9
$ws.Range("K102").Value = 1642.7576
$ws.Range("L102").Value = 103438.9
$ws.Range("M102").Value = -20.75759999999991
$ws.Range("N102").Value = -106682.9

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5108.3335
$ws.Range("I122").Value = 4720
$ws.Range("J122").Value = 5385.7144
$ws.Range("K122").Value = 14160
$ws.Range("L122").Value = 16157.1432
$ws.Range("M122").Value = -11710
$ws.Range("N122").Value = -21057.1432

# Sheet LTW, row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 418666.72
$ws.Range("I2").Value = 538461.5600000001
$ws.Range("J2").Value = 107200.2
$ws.Range("K2").Value = 538461.5600000001
$ws.Range("L2").Value = 107200.2
$ws.Range("M2").Value = -538349.5600000001
$ws.Range("N2").Value = -107424.2

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2945.2778
$ws.Range("I122").Value = 2421.6667
$ws.Range("J122").Value = 3992.5
$ws.Range("K122").Value = 7265.000100000001
$ws.Range("L122").Value = 11977.5
$ws.Range("M122").Value = -4815.000100000001
$ws.Range("N122").Value = -16877.5

# Sheet WVR, row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 94098750
$ws.Range("J2").Value = 3583572.5
$ws.Range("L2").Value = 3583572.5
$ws.Range("N2").Value = -3583796.5

# Sheet WVR, row 116
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 40000
$ws.Range("J116").Value = 40000
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 391255.22
$ws.Range("I132").Value = 1251983
$ws.Range("K132").Value = 3755949
$ws.Range("M132").Value = -3753419
